# Move the title text ("Ridge Regresstion" / "Lasso Regression") from
# row 1 (A1 / E1) down into row 2 (A2 / E2), matching the formatting of
# the header row it now shares (B2/C2/F2/G2), leaving A1 / E1 blank but
# keeping their original cell formatting. Also widen columns A and E to
# fit the relocated text, and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the title text currently sitting in A1 / E1 (.Value2 avoids
# the quirky boxed-Variant echo that plain .Value can produce on read).
$leftTitle = $ws.Range("A1").Value2
$rightTitle = $ws.Range("E1").Value2

# Clear the text from A1 / E1 but keep their existing formatting/style.
$ws.Range("A1").Value = $null
$ws.Range("E1").Value = $null

# Give the new A2 / E2 cells the same formatting as the rest of row 2
# (copy format only, so the shared style is reused rather than a new one
# being created) before filling in the relocated title text.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

$ws.Range("A2").Value = $leftTitle
$ws.Range("E2").Value = $rightTitle

# Widen columns A and E so the relocated titles fit nicely.
$ws.Columns.Item(1).ColumnWidth = 20.833333333333332
$ws.Columns.Item(5).ColumnWidth = 18.833333333333332

# Update the active selection to match (I11 -> I10).
$ws.Range("I10").Select()
